$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 166032
$ws.Range("C4").Value = 156953
$ws.Range("C5").Value = 9080
$ws.Range("C8").Value = 65.06999999999999
